$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-23: the placeholder "date-looking" text in column A is replaced
#     with the real row labels (titles / footnotes / period labels) ---
$ws.Range("A2").Value = '                National Gem and Jewellery Authority'
$ws.Range("A3").Value = '               Sri Lanka Customs'
$ws.Range("A4").Value = '              Central Bank of Sri Lanka'
$ws.Range("A5").Value = '(a)  The latest version of SITC Revision 4 published in 2006'
$ws.Range("A6").Value = '(b)  Provisional'
$ws.Range("A7").Value = '2.04: Import Performance based on the Standard International Trade Classification (SITC) Monthly 2014-2024 (a)'
$ws.Range("A8").Value = '2.04: Import Performance based on the Standard International Trade Classification (SITC) Monthly 2014-2024 (a)'

# "2023"/"2024" must stay text (not auto-convert to numbers): force text
# format, assign, then drop back to the Normal style so no stray number
# format sticks around on the cell.
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '2023'
$ws.Range("A9").Style = "Normal"

$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '2024'
$ws.Range("A10").Style = "Normal"

$ws.Range("A11").Value = 'January'
$ws.Range("A12").Value = 'January'
$ws.Range("A13").Value = 'January'
$ws.Range("A14").Value = 'January'
$ws.Range("A15").Value = 'January'
$ws.Range("A16").Value = 'January'
$ws.Range("A17").Value = 'January'
$ws.Range("A18").Value = 'January'
$ws.Range("A19").Value = 'January'
$ws.Range("A20").Value = 'January'
$ws.Range("A21").Value = 'January'
$ws.Range("A22").Value = 'Sources: Ceylon Petroleum Corporation and Other Exporters of Petroleum'
$ws.Range("A23").Value = 'Table 2.04.4: Imports (Rs Million)'

# --- Rows 24-141: the leftover placeholder text in column A is removed
#     outright (the numeric data columns D.. are left untouched) ---
$ws.Range("A24:A141").ClearContents()

# A handful of rows in that range never had any numeric data either (they
# were blank filler rows between month groups), so clearing column A leaves
# them completely empty. Touch a default (no-op) format on each one so it
# stays registered as a real (blank) row in the sheet instead of being
# trimmed away entirely.
$ws.Range("A27").Font.Bold = $false
$ws.Range("A40").Font.Bold = $false
$ws.Range("A53").Font.Bold = $false
$ws.Range("A66").Font.Bold = $false
$ws.Range("A79").Font.Bold = $false
$ws.Range("A92").Font.Bold = $false
$ws.Range("A105").Font.Bold = $false
$ws.Range("A118").Font.Bold = $false
$ws.Range("A131").Font.Bold = $false

# --- Rows 142-159: these rows only ever carried leftover column-A text and
#     no numeric data, so once column A is cleared they become fully blank
#     rows. Touch the same no-op format so the rows stay registered as
#     part of the sheet's used range instead of being trimmed away. ---
$ws.Range("A142:A159").ClearContents()
$ws.Range("A142").Font.Bold = $false
$ws.Range("A143").Font.Bold = $false
$ws.Range("A144").Font.Bold = $false
$ws.Range("A145").Font.Bold = $false
$ws.Range("A146").Font.Bold = $false
$ws.Range("A147").Font.Bold = $false
$ws.Range("A148").Font.Bold = $false
$ws.Range("A149").Font.Bold = $false
$ws.Range("A150").Font.Bold = $false
$ws.Range("A151").Font.Bold = $false
$ws.Range("A152").Font.Bold = $false
$ws.Range("A153").Font.Bold = $false
$ws.Range("A154").Font.Bold = $false
$ws.Range("A155").Font.Bold = $false
$ws.Range("A156").Font.Bold = $false
$ws.Range("A157").Font.Bold = $false
$ws.Range("A158").Font.Bold = $false
$ws.Range("A159").Font.Bold = $false

# --- Rows 160-325: delete these rows entirely, shrinking the sheet's used
#     range down to A1:T159 ---
$ws.Range("A160:A325").EntireRow.Delete()
